$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing id values in rows 2-11: 101-110 -> 501-510
$ws.Range("A2").Value = 501
$ws.Range("A3").Value = 502
$ws.Range("A4").Value = 503
$ws.Range("A5").Value = 504
$ws.Range("A6").Value = 505
$ws.Range("A7").Value = 506
$ws.Range("A8").Value = 507
$ws.Range("A9").Value = 508
$ws.Range("A10").Value = 509
$ws.Range("A11").Value = 510

# Add new rows 12-64 with poker hand / card data
$ws.Range("A12").Value = 101
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "Royal Flush"
$ws.Range("D12").Value = 200
$ws.Range("E12").Value = 12
$ws.Range("A13").Value = 102
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "Straight Flush"
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = 8
$ws.Range("A14").Value = 103
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = "Four of a Kind"
$ws.Range("D14").Value = 60
$ws.Range("E14").Value = 7
$ws.Range("H14").Value = "所有包含7的牌型得分翻倍"
$ws.Range("A15").Value = 104
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = "Hulu"
$ws.Range("D15").Value = 40
$ws.Range("E15").Value = 4
$ws.Range("H15").Value = "所有的对牌（Pair）得分翻倍"
$ws.Range("A16").Value = 105
$ws.Range("B16").Value = 5
$ws.Range("C16").Value = "SAmehua"
$ws.Range("D16").Value = 35
$ws.Range("E16").Value = 4
$ws.Range("H16").Value = "所有同花顺（Flush）的牌型得分增加50%"
$ws.Range("A17").Value = 106
$ws.Range("B17").Value = 6
$ws.Range("C17").Value = "Flush"
$ws.Range("D17").Value = 30
$ws.Range("E17").Value = 4
$ws.Range("H17").Value = "所有顺子（Straight）的牌型得分翻倍"
$ws.Range("A18").Value = 107
$ws.Range("B18").Value = 7
$ws.Range("C18").Value = "Three of a kind"
$ws.Range("D18").Value = 30
$ws.Range("E18").Value = 3
$ws.Range("H18").Value = "每个玩家的最终得分随机增加1.5到3倍"
$ws.Range("A19").Value = 108
$ws.Range("B19").Value = 8
$ws.Range("C19").Value = "Two Pairs"
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = 2
$ws.Range("H19").Value = "包含至少一张红心牌的牌型得分增加50%"
$ws.Range("A20").Value = 109
$ws.Range("B20").Value = 9
$ws.Range("C20").Value = "Pairs"
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = 2
$ws.Range("H20").Value = "任何包含A的牌型得分增加50%"
$ws.Range("A21").Value = 110
$ws.Range("B21").Value = 10
$ws.Range("C21").Value = "High Card"
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 1
$ws.Range("H21").Value = "所有四张同牌（Four of a Kind）的牌型得分增加两倍"
$ws.Range("A22").Value = 211
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = "HT"
$ws.Range("D22").Value = 10
$ws.Range("A23").Value = 212
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = "HJ"
$ws.Range("D23").Value = 11
$ws.Range("A24").Value = 213
$ws.Range("B24").Value = 2
$ws.Range("C24").Value = "HQ"
$ws.Range("D24").Value = 12
$ws.Range("A25").Value = 214
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = "HK"
$ws.Range("D25").Value = 13
$ws.Range("A26").Value = 215
$ws.Range("B26").Value = 2
$ws.Range("C26").Value = "HA"
$ws.Range("D26").Value = 14
$ws.Range("A27").Value = 216
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = "S2"
$ws.Range("D27").Value = 2
$ws.Range("A28").Value = 217
$ws.Range("B28").Value = 2
$ws.Range("C28").Value = "S3"
$ws.Range("D28").Value = 3
$ws.Range("A29").Value = 218
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = "S4"
$ws.Range("D29").Value = 4
$ws.Range("A30").Value = 219
$ws.Range("B30").Value = 2
$ws.Range("C30").Value = "S5"
$ws.Range("D30").Value = 5
$ws.Range("A31").Value = 220
$ws.Range("B31").Value = 2
$ws.Range("C31").Value = "S6"
$ws.Range("D31").Value = 6
$ws.Range("A32").Value = 221
$ws.Range("B32").Value = 2
$ws.Range("C32").Value = "S7"
$ws.Range("D32").Value = 7
$ws.Range("A33").Value = 222
$ws.Range("B33").Value = 2
$ws.Range("C33").Value = "S8"
$ws.Range("D33").Value = 8
$ws.Range("A34").Value = 223
$ws.Range("B34").Value = 2
$ws.Range("C34").Value = "S9"
$ws.Range("D34").Value = 9
$ws.Range("A35").Value = 224
$ws.Range("B35").Value = 2
$ws.Range("C35").Value = "ST"
$ws.Range("D35").Value = 10
$ws.Range("A36").Value = 225
$ws.Range("B36").Value = 2
$ws.Range("C36").Value = "SJ"
$ws.Range("D36").Value = 11
$ws.Range("A37").Value = 226
$ws.Range("B37").Value = 2
$ws.Range("C37").Value = "SQ"
$ws.Range("D37").Value = 12
$ws.Range("A38").Value = 227
$ws.Range("B38").Value = 2
$ws.Range("C38").Value = "SA"
$ws.Range("D38").Value = 13
$ws.Range("A39").Value = 228
$ws.Range("B39").Value = 2
$ws.Range("C39").Value = "D2"
$ws.Range("D39").Value = 2
$ws.Range("A40").Value = 229
$ws.Range("B40").Value = 2
$ws.Range("C40").Value = "D3"
$ws.Range("D40").Value = 3
$ws.Range("A41").Value = 230
$ws.Range("B41").Value = 2
$ws.Range("C41").Value = "D4"
$ws.Range("D41").Value = 4
$ws.Range("A42").Value = 231
$ws.Range("B42").Value = 2
$ws.Range("C42").Value = "D5"
$ws.Range("D42").Value = 5
$ws.Range("A43").Value = 232
$ws.Range("B43").Value = 2
$ws.Range("C43").Value = "D6"
$ws.Range("D43").Value = 6
$ws.Range("A44").Value = 233
$ws.Range("B44").Value = 2
$ws.Range("C44").Value = "D7"
$ws.Range("D44").Value = 7
$ws.Range("A45").Value = 234
$ws.Range("B45").Value = 2
$ws.Range("C45").Value = "D8"
$ws.Range("D45").Value = 8
$ws.Range("A46").Value = 235
$ws.Range("B46").Value = 2
$ws.Range("C46").Value = "D9"
$ws.Range("D46").Value = 9
$ws.Range("A47").Value = 236
$ws.Range("B47").Value = 2
$ws.Range("C47").Value = "DT"
$ws.Range("D47").Value = 10
$ws.Range("A48").Value = 237
$ws.Range("B48").Value = 2
$ws.Range("C48").Value = "DJ"
$ws.Range("D48").Value = 11
$ws.Range("A49").Value = 238
$ws.Range("B49").Value = 2
$ws.Range("C49").Value = "DQ"
$ws.Range("D49").Value = 12
$ws.Range("A50").Value = 239
$ws.Range("B50").Value = 2
$ws.Range("C50").Value = "DK"
$ws.Range("D50").Value = 13
$ws.Range("A51").Value = 240
$ws.Range("B51").Value = 2
$ws.Range("C51").Value = "DA"
$ws.Range("D51").Value = 14
$ws.Range("A52").Value = 241
$ws.Range("B52").Value = 2
$ws.Range("C52").Value = "C2"
$ws.Range("D52").Value = 2
$ws.Range("A53").Value = 242
$ws.Range("B53").Value = 2
$ws.Range("C53").Value = "C3"
$ws.Range("D53").Value = 3
$ws.Range("A54").Value = 243
$ws.Range("B54").Value = 2
$ws.Range("C54").Value = "C4"
$ws.Range("D54").Value = 4
$ws.Range("A55").Value = 244
$ws.Range("B55").Value = 2
$ws.Range("C55").Value = "C5"
$ws.Range("D55").Value = 5
$ws.Range("A56").Value = 245
$ws.Range("B56").Value = 2
$ws.Range("C56").Value = "C6"
$ws.Range("D56").Value = 6
$ws.Range("A57").Value = 246
$ws.Range("B57").Value = 2
$ws.Range("C57").Value = "C7"
$ws.Range("D57").Value = 7
$ws.Range("A58").Value = 247
$ws.Range("B58").Value = 2
$ws.Range("C58").Value = "C8"
$ws.Range("D58").Value = 8
$ws.Range("A59").Value = 248
$ws.Range("B59").Value = 2
$ws.Range("C59").Value = "C9"
$ws.Range("D59").Value = 9
$ws.Range("A60").Value = 249
$ws.Range("B60").Value = 2
$ws.Range("C60").Value = "CT"
$ws.Range("D60").Value = 10
$ws.Range("A61").Value = 250
$ws.Range("B61").Value = 2
$ws.Range("C61").Value = "CJ"
$ws.Range("D61").Value = 11
$ws.Range("A62").Value = 251
$ws.Range("B62").Value = 2
$ws.Range("C62").Value = "CQ"
$ws.Range("D62").Value = 12
$ws.Range("A63").Value = 252
$ws.Range("B63").Value = 2
$ws.Range("C63").Value = "CK"
$ws.Range("D63").Value = 13
$ws.Range("A64").Value = 253
$ws.Range("B64").Value = 2
$ws.Range("C64").Value = "CA"
$ws.Range("D64").Value = 14

# Update the active selection to match the final state
$ws.Range("D21").Select()
